# Update "Update countries & provincias Spain"
# - Refresh the "last updated" timestamp.
# - Refresh the case counts for Asturias and Murcia; since the table is kept
#   sorted by "Casos totales" (column B) descending, both provinces move up
#   in rank, which shifts the province names/labels shown on the rows around
#   them while the rows that held those labels now show the displaced names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 22:29"

# --- Asturias moves from rank 20 (row 20) up to rank 18 (row 18) -------
# New data for Asturias:
$ws.Range("A18").Value = "Asturias"
$ws.Range("B18").Value = 1088
$ws.Range("C18").Value = 76
$ws.Range("D18").Value = 971
$ws.Range("E18").Value = 41

# Malaga shifts down one rank (row 18 -> row 19), values unchanged:
$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 1053
$ws.Range("C19").Value = 80
$ws.Range("D19").Value = 917
$ws.Range("E19").Value = 56

# Tenerife shifts down one rank (row 19 -> row 20), values unchanged:
$ws.Range("A20").Value = "Tenerife"
$ws.Range("B20").Value = 1025
$ws.Range("C20").Value = 25
$ws.Range("D20").Value = 964
$ws.Range("E20").Value = 36

# Row 21 (Pontevedra) is unaffected.

# --- Murcia moves from rank 30 (row 30) up to rank 27 (row 27) ---------
# New data for Murcia:
$ws.Range("A27").Value = "Murcia"
$ws.Range("B27").Value = 836
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 804
$ws.Range("E27").Value = 20

# Sevilla shifts down one rank (row 27 -> row 28), values unchanged:
$ws.Range("A28").Value = "Sevilla"
$ws.Range("B28").Value = 830
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 791
$ws.Range("E28").Value = 26

# Valladolid shifts down one rank (row 28 -> row 29), values unchanged:
$ws.Range("A29").Value = "Valladolid"
$ws.Range("B29").Value = 807
$ws.Range("C29").Value = 114
$ws.Range("D29").Value = 648
$ws.Range("E29").Value = 45

# Granada shifts down one rank (row 29 -> row 30), values unchanged:
$ws.Range("A30").Value = "Granada"
$ws.Range("B30").Value = 806
$ws.Range("C30").Value = 11
$ws.Range("D30").Value = 746
$ws.Range("E30").Value = 49

# Row 31 (Leon) is unaffected.

$wb.Save()
